# Automatische test-sync: 2025-08-13 22:42:50
# Appends the newest "Demo inplannen" log entry to the Logs sheet,
# extends the conditional-formatting ranges to cover the new row,
# and bumps the Dashboard summary count for the matching category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 28

$ws.Cells.Item($newRow, 1).Value = "Demo inplannen"
$ws.Cells.Item($newRow, 2).Value = "klantenservice@testbedrijf123.nl"
$ws.Cells.Item($newRow, 3).Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$ws.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Cells.Item($newRow, 6).Value = "2025-08-13 22:42:46"
$ws.Cells.Item($newRow, 7).Value = "Nee"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional formatting ranges (D,G,H,I,J) from row 27 to row 28
# so the newly appended row keeps the same highlighting rules.
$dCond = $ws.Range("D2:D27").FormatConditions
$dCond.Item(1).ModifyAppliesToRange($ws.Range("D2:D28"))

$gCond = $ws.Range("G2:G27").FormatConditions
$gCond.Item(1).ModifyAppliesToRange($ws.Range("G2:G28"))

$hCond = $ws.Range("H2:H27").FormatConditions
$hCond.Item(1).ModifyAppliesToRange($ws.Range("H2:H28"))

$iCond = $ws.Range("I2:I27").FormatConditions
$iCond.Item(1).ModifyAppliesToRange($ws.Range("I2:I28"))

$jCond = $ws.Range("J2:J27").FormatConditions
$jCond.Item(1).ModifyAppliesToRange($ws.Range("J2:J28"))

# Update the Dashboard summary count for "Intern verzoek / Actie voor medewerker"
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Cells.Item(2, 2).Value = 27
